$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: "...iremos por em discussão..." -> split into three runs:
#   "rimeiro, iremos " | "colocar" | " em discussão: ..."
# (word "por" replaced with "colocar")
# ---------------------------------------------------------------
$full = $d.Content.Text
$offset1 = $full.IndexOf("por em discuss")
$rWord = $d.Range($offset1, $offset1 + 3)
$rWord.Text = "colocar"
# Touch formatting on the just-inserted word so the engine keeps it
# as its own run (mirrors a real edit: select "por", type "colocar").
$rNew = $d.Range($offset1, $offset1 + 7)
$rNew.Bold = 1
$rNew.Bold = 0

# ---------------------------------------------------------------
# Change 2: split the SCRUM paragraph into two runs:
#   "3. Poderia ser aplicado o SCRUM" | " porque é nela que fazemos
#    todo o planejamento detalhado e sabemos a função de cada
#    membro desse projeto."
# ---------------------------------------------------------------
$full = $d.Content.Text
$oldTail = "3. Poderia ser aplicado o SCRUM. Sempre é recomendado em qualquer projeto ter sempre um PO(Product Owner) que seria o responsável pela garantia do investimento e também por estar em contato direto com o cliente pra suprir suas necessidades, um SM(Scrum Master) que seria o responsável pela aprovação do projeto e da tomada de decisão final e por fim, todo o time que será feito para execução do projeto."
$keep = "3. Poderia ser aplicado o SCRUM"
$offset2 = $full.IndexOf($oldTail)
$rTail = $d.Range($offset2 + $keep.Length, $offset2 + $oldTail.Length)
$rTail.Text = " porque é nela que fazemos todo o planejamento detalhado e sabemos a função de cada membro desse projeto."
$rKeep = $d.Range($offset2, $offset2 + $keep.Length)
$rKeep.Bold = 1
$rKeep.Bold = 0

# ---------------------------------------------------------------
# Change 3: the trailing empty paragraph becomes paragraph "4. ..."
# ---------------------------------------------------------------
$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Item($paragraphs.Count)
$lastPara.Range.Text = "4. O time seria separado em: Um PO(Product Owner) e ele será o responsável pela aprovação das ideias e também da comunicação direta com o Scrum Master(SU). O SU(Scrum Master) será o chefe de todo o projeto, é nele que toda tomada final de decisão é decidida para continuação do projeto, todos os ajustes é passado para ele pra verificar e ver o melhor. E por fim tem todo o time: Time pra verificar em tempo real as condições de tráfego, time pra localização e atualização das rotas para os caminhoneiros, time pra determinar a quantidade de lixo que será coletado e as ruas que serão coletadas e um time pra coletar o histórico de coleta e guardar pras futuras coletas."
